# Insert a new weekly observation row for Puerro (Vega Central Mapocho de
# Santiago) above the existing row 63. This pushes the previous rows
# 63-75 down to 64-76 (dimension grows from A1:R75 to A1:R76) and fills
# the freshly inserted row 63 with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("63:63").Insert()

$ws.Range("A63").Value = 9
$ws.Range("B63").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C63").Value = "Metropolitana"
$ws.Range("D63").Value = 44505
$ws.Range("E63").Value = 13
$ws.Range("F63").Value = 100112005
$ws.Range("G63").Value = "Puerro"
$ws.Range("H63").Value = "Sin especificar"
$ws.Range("I63").Value = "Primera"
$ws.Range("J63").Value = 160
$ws.Range("K63").Value = 6000
$ws.Range("L63").Value = 7000
$ws.Range("M63").Value = 6500
$ws.Range("N63").Value = "$/paquete 20 unidades"
$ws.Range("O63").Value = "Provincia de Chacabuco"
$ws.Range("P63").Value = 325
$ws.Range("Q63").Value = 20
$ws.Range("R63").Value = "Hortaliza"
